# Curation WP5166 WP5169 WP5173 WP5190
# Update ALAS -> ALAS2 (rows 2-3), add Uniprot accession + curator note,
# and fill in previously-missing Rhea/CHEBI ids for the UROD substrate rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (EnzymePW "ALAS" renamed to "ALAS2", its official gene symbol)
$ws.Range("A2").Value = "ALAS2"
$ws.Range("D2").Value = "P22557"
$ws.Range("N2").Value = "“ usually only one isoform of bacterial ALAS proteins, mammals carry two different ALAS isoforms, one of them fulfilling a housekeeping function (ALAS1) and the other sustaining high levels of heme biosynthesis in erythrocytes (ALAS2); PMID: 2347585”"

# Row 3 (second kinetic entry for the same enzyme)
$ws.Range("A3").Value = "ALAS2"
$ws.Range("D3").Value = "P22557"
$ws.Range("N3").Value = "“ usually only one isoform of bacterial ALAS proteins, mammals carry two different ALAS isoforms, one of them fulfilling a housekeeping function (ALAS1) and the other sustaining high levels of heme biosynthesis in erythrocytes (ALAS2); PMID: 2347585”"

# Rows 9-16 (UROD substrate entries): fill in the previously blank RheaID / CHEBIID values
$ws.Range("F9").Value = "RHEA:31240"

$ws.Range("F10").Value = "RHEA:31240"
$ws.Range("G10").Value = "CHEBI:89912"

$ws.Range("F11").Value = "RHEA:31240"
$ws.Range("G11").Value = "CHEBI:89912"

$ws.Range("F12").Value = "RHEA:31240"
$ws.Range("G12").Value = "CHEBI:89912"

$ws.Range("F13").Value = "RHEA:31240"
$ws.Range("G13").Value = "CHEBI:89912"

$ws.Range("F14").Value = "RHEA:31240"
$ws.Range("G14").Value = "CHEBI:89912"

$ws.Range("F15").Value = "RHEA:31240"

$ws.Range("F16").Value = "RHEA:31240"

# View state: zoom in and move the active selection to E14
$ws.Application.ActiveWindow.DisplayGridlines = $true
$ws.Application.ActiveWindow.Zoom = 95
$ws.Range("E14").Select()
